# Update "Platform Coverage" sheet (sheet1) from a yearly (2018-2040, cols H:AD)
# coverage timeline to a 6-monthly timeline (2018-2040 step 0.5, cols H:AZ),
# per commit "6 monthly coverage scenario 2 and 3a".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Row 1: header years, H1:AZ1, starting at 2018 and incrementing by 0.5.
# (was H1:AD1, 2018..2040 by whole years)
# ---------------------------------------------------------------------------
$year = 2018.0
for ($c = 8; $c -le 52; $c++) {
    if ($year -eq [Math]::Floor($year)) {
        $ws.Cells.Item(1, $c).Value = [int]$year
    } else {
        $ws.Cells.Item(1, $c).Value = $year
    }
    $year += 0.5
}

# ---------------------------------------------------------------------------
# Row 2: coverage 0.6 now also continues through every whole-year column
# up to column V (2025) -- i.e. add P2, R2, T2, V2 = 0.6 (H2,J2,L2,N2
# already held 0.6 and are untouched).
# ---------------------------------------------------------------------------
foreach ($c in 16,18,20,22) {
    $ws.Cells.Item(2, $c).Value = 0.6
}

# ---------------------------------------------------------------------------
# Rows 3, 4, 5: coverage that used to be recorded only on whole-year columns
# P,R,T,V,X,Z,AB,AD now needs to be recorded on every half-year column from
# X (2026) through AZ (2040) instead. Clear the old whole-year cells P..V,
# then fill X..AZ (29 columns) with the row's coverage value.
# ---------------------------------------------------------------------------
foreach ($c in 16,18,20,22) {
    $ws.Cells.Item(3, $c).ClearContents()
    $ws.Cells.Item(4, $c).ClearContents()
    $ws.Cells.Item(5, $c).ClearContents()
}

for ($c = 24; $c -le 52; $c++) {
    $ws.Cells.Item(3, $c).Value = 0.75
    $ws.Cells.Item(4, $c).Value = 0.5
    $ws.Cells.Item(5, $c).Value = 0.5
}

# ---------------------------------------------------------------------------
# Update the recorded selection to match the authored workbook state.
# ---------------------------------------------------------------------------
$ws.Range("J14").Select() | Out-Null
